$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 4.979788333333333
$ws.Cells.Item(2, 8).Value = 14.939365
$ws.Cells.Item(2, 9).Value = 0.129176854764059
$ws.Cells.Item(2, 10).Value = 0.129176854764059
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 161.7750676666667
$ws.Cells.Item(2, 14).Value = 485.325203
$ws.Cells.Item(2, 15).Value = 0.9790864123038654
$ws.Cells.Item(2, 16).Value = 0.9790864123038654
$ws.Cells.Item(2, 17).Value = 805.6055945906772
$ws.Cells.Item(2, 18).Value = 7250.450351316095
$ws.Cells.Item(2, 19).Value = 0.12647530328364
$ws.Cells.Item(2, 20).Value = 0.12647530328364

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 4.979788333333333
$ws.Cells.Item(3, 8).Value = 14.939365
$ws.Cells.Item(3, 9).Value = 0.129176854764059
$ws.Cells.Item(3, 10).Value = 0.129176854764059
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 0.67894
$ws.Cells.Item(3, 14).Value = 2.03682
$ws.Cells.Item(3, 15).Value = 0.004109044356199978
$ws.Cells.Item(3, 16).Value = 0.004109044356199979
$ws.Cells.Item(3, 17).Value = 3.380977491033333
$ws.Cells.Item(3, 18).Value = 30.4287974193
$ws.Cells.Item(3, 19).Value = 0.0005307934260199207
$ws.Cells.Item(3, 20).Value = 0.000530793426019921

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 4.979788333333333
$ws.Cells.Item(4, 8).Value = 14.939365
$ws.Cells.Item(4, 9).Value = 0.129176854764059
$ws.Cells.Item(4, 10).Value = 0.129176854764059
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 1.763201333333333
$ws.Cells.Item(4, 14).Value = 5.289604
$ws.Cells.Item(4, 15).Value = 0.01067115280816804
$ws.Cells.Item(4, 16).Value = 0.01067115280816804
$ws.Cells.Item(4, 17).Value = 8.780369429051111
$ws.Cells.Item(4, 18).Value = 79.02332486146
$ws.Cells.Item(4, 19).Value = 0.001378465956465803
$ws.Cells.Item(4, 20).Value = 0.001378465956465803

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 4.979788333333333
$ws.Cells.Item(5, 8).Value = 14.939365
$ws.Cells.Item(5, 9).Value = 0.129176854764059
$ws.Cells.Item(5, 10).Value = 0.129176854764059
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 1.013424
$ws.Cells.Item(5, 14).Value = 3.040272
$ws.Cells.Item(5, 15).Value = 0.006133390531766587
$ws.Cells.Item(5, 16).Value = 0.006133390531766588
$ws.Cells.Item(5, 17).Value = 5.04663701192
$ws.Cells.Item(5, 18).Value = 45.41973310728
$ws.Cells.Item(5, 19).Value = 0.0007922920979332669
$ws.Cells.Item(5, 20).Value = 0.0007922920979332671

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 18.019504
$ws.Cells.Item(6, 8).Value = 54.058512
$ws.Cells.Item(6, 9).Value = 0.467430078412646
$ws.Cells.Item(6, 10).Value = 0.4674300784126461
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 161.7750676666667
$ws.Cells.Item(6, 14).Value = 485.325203
$ws.Cells.Item(6, 15).Value = 0.9790864123038654
$ws.Cells.Item(6, 16).Value = 0.9790864123038654
$ws.Cells.Item(6, 17).Value = 2915.106478919771
$ws.Cells.Item(6, 18).Value = 26235.95831027794
$ws.Cells.Item(6, 19).Value = 0.4576544384759521
$ws.Cells.Item(6, 20).Value = 0.4576544384759521

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 18.019504
$ws.Cells.Item(7, 8).Value = 54.058512
$ws.Cells.Item(7, 9).Value = 0.467430078412646
$ws.Cells.Item(7, 10).Value = 0.4674300784126461
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 0.67894
$ws.Cells.Item(7, 14).Value = 2.03682
$ws.Cells.Item(7, 15).Value = 0.004109044356199978
$ws.Cells.Item(7, 16).Value = 0.004109044356199979
$ws.Cells.Item(7, 17).Value = 12.23416204576
$ws.Cells.Item(7, 18).Value = 110.10745841184
$ws.Cells.Item(7, 19).Value = 0.001920690925619596
$ws.Cells.Item(7, 20).Value = 0.001920690925619597

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 18.019504
$ws.Cells.Item(8, 8).Value = 54.058512
$ws.Cells.Item(8, 9).Value = 0.467430078412646
$ws.Cells.Item(8, 10).Value = 0.4674300784126461
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 1.763201333333333
$ws.Cells.Item(8, 14).Value = 5.289604
$ws.Cells.Item(8, 15).Value = 0.01067115280816804
$ws.Cells.Item(8, 16).Value = 0.01067115280816804
$ws.Cells.Item(8, 17).Value = 31.77201347880533
$ws.Cells.Item(8, 18).Value = 285.948121309248
$ws.Cells.Item(8, 19).Value = 0.004988017793875315
$ws.Cells.Item(8, 20).Value = 0.004988017793875316

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 18.019504
$ws.Cells.Item(9, 8).Value = 54.058512
$ws.Cells.Item(9, 9).Value = 0.467430078412646
$ws.Cells.Item(9, 10).Value = 0.4674300784126461
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 1.013424
$ws.Cells.Item(9, 14).Value = 3.040272
$ws.Cells.Item(9, 15).Value = 0.006133390531766587
$ws.Cells.Item(9, 16).Value = 0.006133390531766588
$ws.Cells.Item(9, 17).Value = 18.261397821696
$ws.Cells.Item(9, 18).Value = 164.352580395264
$ws.Cells.Item(9, 19).Value = 0.002866931217199037
$ws.Cells.Item(9, 20).Value = 0.002866931217199037

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 8.752692000000001
$ws.Cells.Item(10, 8).Value = 26.258076
$ws.Cells.Item(10, 9).Value = 0.2270468436801446
$ws.Cells.Item(10, 10).Value = 0.2270468436801446
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 161.7750676666667
$ws.Cells.Item(10, 14).Value = 485.325203
$ws.Cells.Item(10, 15).Value = 0.9790864123038654
$ws.Cells.Item(10, 16).Value = 0.9790864123038654
$ws.Cells.Item(10, 17).Value = 1415.967340565492
$ws.Cells.Item(10, 18).Value = 12743.70606508943
$ws.Cells.Item(10, 19).Value = 0.2222984796037093
$ws.Cells.Item(10, 20).Value = 0.2222984796037093

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 8.752692000000001
$ws.Cells.Item(11, 8).Value = 26.258076
$ws.Cells.Item(11, 9).Value = 0.2270468436801446
$ws.Cells.Item(11, 10).Value = 0.2270468436801446
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 0.67894
$ws.Cells.Item(11, 14).Value = 2.03682
$ws.Cells.Item(11, 15).Value = 0.004109044356199978
$ws.Cells.Item(11, 16).Value = 0.004109044356199979
$ws.Cells.Item(11, 17).Value = 5.942552706480001
$ws.Cells.Item(11, 18).Value = 53.48297435832001
$ws.Cells.Item(11, 19).Value = 0.0009329455516169168
$ws.Cells.Item(11, 20).Value = 0.000932945551616917

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 8.752692000000001
$ws.Cells.Item(12, 8).Value = 26.258076
$ws.Cells.Item(12, 9).Value = 0.2270468436801446
$ws.Cells.Item(12, 10).Value = 0.2270468436801446
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 1.763201333333333
$ws.Cells.Item(12, 14).Value = 5.289604
$ws.Cells.Item(12, 15).Value = 0.01067115280816804
$ws.Cells.Item(12, 16).Value = 0.01067115280816804
$ws.Cells.Item(12, 17).Value = 15.432758204656
$ws.Cells.Item(12, 18).Value = 138.894823841904
$ws.Cells.Item(12, 19).Value = 0.002422851563523065
$ws.Cells.Item(12, 20).Value = 0.002422851563523066

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 8.752692000000001
$ws.Cells.Item(13, 8).Value = 26.258076
$ws.Cells.Item(13, 9).Value = 0.2270468436801446
$ws.Cells.Item(13, 10).Value = 0.2270468436801446
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 1.013424
$ws.Cells.Item(13, 14).Value = 3.040272
$ws.Cells.Item(13, 15).Value = 0.006133390531766587
$ws.Cells.Item(13, 16).Value = 0.006133390531766588
$ws.Cells.Item(13, 17).Value = 8.870188137408
$ws.Cells.Item(13, 18).Value = 79.83169323667201
$ws.Cells.Item(13, 19).Value = 0.001392566961295287
$ws.Cells.Item(13, 20).Value = 0.001392566961295287

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 6.798175000000001
$ws.Cells.Item(14, 8).Value = 20.394525
$ws.Cells.Item(14, 9).Value = 0.1763462231431503
$ws.Cells.Item(14, 10).Value = 0.1763462231431503
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 161.7750676666667
$ws.Cells.Item(14, 14).Value = 485.325203
$ws.Cells.Item(14, 15).Value = 0.9790864123038654
$ws.Cells.Item(14, 16).Value = 0.9790864123038654
$ws.Cells.Item(14, 17).Value = 1099.775220634842
$ws.Cells.Item(14, 18).Value = 9897.976985713576
$ws.Cells.Item(14, 19).Value = 0.1726581909405639
$ws.Cells.Item(14, 20).Value = 0.1726581909405639

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 6.798175000000001
$ws.Cells.Item(15, 8).Value = 20.394525
$ws.Cells.Item(15, 9).Value = 0.1763462231431503
$ws.Cells.Item(15, 10).Value = 0.1763462231431503
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 0.67894
$ws.Cells.Item(15, 14).Value = 2.03682
$ws.Cells.Item(15, 15).Value = 0.004109044356199978
$ws.Cells.Item(15, 16).Value = 0.004109044356199979
$ws.Cells.Item(15, 17).Value = 4.6155529345
$ws.Cells.Item(15, 18).Value = 41.5399764105
$ws.Cells.Item(15, 19).Value = 0.0007246144529435438
$ws.Cells.Item(15, 20).Value = 0.0007246144529435442

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 6.798175000000001
$ws.Cells.Item(16, 8).Value = 20.394525
$ws.Cells.Item(16, 9).Value = 0.1763462231431503
$ws.Cells.Item(16, 10).Value = 0.1763462231431503
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 1.763201333333333
$ws.Cells.Item(16, 14).Value = 5.289604
$ws.Cells.Item(16, 15).Value = 0.01067115280816804
$ws.Cells.Item(16, 16).Value = 0.01067115280816804
$ws.Cells.Item(16, 17).Value = 11.98655122423333
$ws.Cells.Item(16, 18).Value = 107.8789610181
$ws.Cells.Item(16, 19).Value = 0.001881817494303857
$ws.Cells.Item(16, 20).Value = 0.001881817494303857

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 6.798175000000001
$ws.Cells.Item(17, 8).Value = 20.394525
$ws.Cells.Item(17, 9).Value = 0.1763462231431503
$ws.Cells.Item(17, 10).Value = 0.1763462231431503
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 1.013424
$ws.Cells.Item(17, 14).Value = 3.040272
$ws.Cells.Item(17, 15).Value = 0.006133390531766587
$ws.Cells.Item(17, 16).Value = 0.006133390531766588
$ws.Cells.Item(17, 17).Value = 6.8894337012
$ws.Cells.Item(17, 18).Value = 62.0049033108
$ws.Cells.Item(17, 19).Value = 0.001081600255338996
$ws.Cells.Item(17, 20).Value = 0.001081600255338996
